# Auto-generated script to apply scheduled market-data refresh to Sheets
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 38
$ws.Range("I11").Value = 38
$ws.Range("K11").Value = 38
$ws.Range("M11").Value = 102
$ws.Range("H74").Value = 5338.3687
$ws.Range("I74").Value = 4299.2856
$ws.Range("K74").Value = 4299.2856
$ws.Range("M74").Value = -3363.2856
$ws.Range("H77").Value = 5338.3687
$ws.Range("I77").Value = 4299.2856
$ws.Range("K77").Value = 21496.428
$ws.Range("M77").Value = -16816.428
$ws.Range("H97").Value = 2981
$ws.Range("J97").Value = 2981
$ws.Range("L97").Value = 8943
$ws.Range("N97").Value = -9935
$ws.Range("H111").Value = 2256
$ws.Range("I111").Value = 2710.875
$ws.Range("K111").Value = 8132.625
$ws.Range("M111").Value = -5065.625
$ws.Range("H113").Value = 7624.7144
$ws.Range("I113").Value = 5848.75
$ws.Range("J113").Value = 9992.666999999999
$ws.Range("K113").Value = 5848.75
$ws.Range("L113").Value = 9992.666999999999
$ws.Range("M113").Value = -2594.75
$ws.Range("N113").Value = -16500.667
$ws.Range("H116").Value = 9549.666999999999
$ws.Range("I116").Value = 7989.3335
$ws.Range("J116").Value = 11110
$ws.Range("K116").Value = 7989.3335
$ws.Range("L116").Value = 11110
$ws.Range("M116").Value = -4547.3335
$ws.Range("N116").Value = -17994
$ws.Range("H132").Value = 4793.484
$ws.Range("I132").Value = 4899.9653
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 14699.8959
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -12169.8959
$ws.Range("N132").Value = -14808.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 39747
$ws.Range("I44").Value = 39999
$ws.Range("K44").Value = 39999
$ws.Range("M44").Value = -39511
$ws.Range("H45").Value = 3508.5557
$ws.Range("I45").Value = 3579.4
$ws.Range("K45").Value = 3579.4
$ws.Range("M45").Value = -3202.4
$ws.Range("H48").Value = 269999
$ws.Range("J48").Value = 269999
$ws.Range("L48").Value = 269999
$ws.Range("N48").Value = -270767
$ws.Range("H61").Value = 2559.4443
$ws.Range("I61").Value = 2076.4285
$ws.Range("K61").Value = 2076.4285
$ws.Range("M61").Value = -1864.4285
$ws.Range("H136").Value = 2559.4443
$ws.Range("I136").Value = 2076.4285
$ws.Range("K136").Value = 6229.2855
$ws.Range("M136").Value = -3679.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4569.6
$ws.Range("I20").Value = 2424.5
$ws.Range("J20").Value = 5999.6665
$ws.Range("K20").Value = 2424.5
$ws.Range("L20").Value = 5999.6665
$ws.Range("M20").Value = -2177.5
$ws.Range("N20").Value = -6493.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2284
$ws.Range("I8").Value = 460
$ws.Range("J8").Value = 3500
$ws.Range("K8").Value = 460
$ws.Range("L8").Value = 3500
$ws.Range("M8").Value = -320
$ws.Range("N8").Value = -3780
$ws.Range("H52").Value = 86997.5
$ws.Range("J52").Value = 86997.5
$ws.Range("L52").Value = 86997.5
$ws.Range("N52").Value = -87585.5
$ws.Range("H122").Value = 10528385
$ws.Range("I122").Value = 2587.5
$ws.Range("K122").Value = 7762.5
$ws.Range("M122").Value = -5312.5
$ws.Range("H134").Value = 3173
$ws.Range("I134").Value = 2682.1904
$ws.Range("J134").Value = 5749.75
$ws.Range("K134").Value = 8046.5712
$ws.Range("L134").Value = 17249.25
$ws.Range("M134").Value = -5511.5712
$ws.Range("N134").Value = -22319.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3300
$ws.Range("I5").Value = 375
$ws.Range("K5").Value = 1125
$ws.Range("M5").Value = -1013
$ws.Range("H37").Value = 138158.4
$ws.Range("J37").Value = 138158.4
$ws.Range("L37").Value = 414475.2
$ws.Range("N37").Value = -414699.2
$ws.Range("H39").Value = 1391.6842
$ws.Range("J39").Value = 3250.75
$ws.Range("L39").Value = 9752.25
$ws.Range("N39").Value = -10340.25
$ws.Range("H55").Value = 3001233
$ws.Range("J55").Value = 1849.5
$ws.Range("L55").Value = 5548.5
$ws.Range("N55").Value = -5902.5
$ws.Range("H104").Value = 6953.222
$ws.Range("J104").Value = 6953.222
$ws.Range("L104").Value = 20859.666
$ws.Range("N104").Value = -26101.666
$ws.Range("H113").Value = 1049.7142
$ws.Range("J113").Value = 2145.3333
$ws.Range("L113").Value = 6435.999899999999
$ws.Range("N113").Value = -10775.9999
$ws.Range("H135").Value = 3300
$ws.Range("I135").Value = 375
$ws.Range("K135").Value = 3375
$ws.Range("M135").Value = -840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8703.375
$ws.Range("I70").Value = 4813.3335
$ws.Range("J70").Value = 11037.4
$ws.Range("K70").Value = 4813.3335
$ws.Range("L70").Value = 11037.4
$ws.Range("M70").Value = -4543.3335
$ws.Range("N70").Value = -11577.4
$ws.Range("H73").Value = 8703.375
$ws.Range("I73").Value = 4813.3335
$ws.Range("J73").Value = 11037.4
$ws.Range("K73").Value = 4813.3335
$ws.Range("L73").Value = 11037.4
$ws.Range("M73").Value = -3877.3335
$ws.Range("N73").Value = -12909.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6774.385
$ws.Range("I7").Value = 6507.909
$ws.Range("K7").Value = 6507.909
$ws.Range("M7").Value = -6395.909
$ws.Range("H16").Value = 848.4286
$ws.Range("I16").Value = 830.24
$ws.Range("K16").Value = 830.24
$ws.Range("M16").Value = -660.24
$ws.Range("H43").Value = 20025
$ws.Range("J43").Value = 20537.5
$ws.Range("L43").Value = 20537.5
$ws.Range("N43").Value = -20923.5
$ws.Range("H82").Value = 2233.4285
$ws.Range("I82").Value = 1529.4
$ws.Range("J82").Value = 2624.5557
$ws.Range("K82").Value = 1529.4
$ws.Range("L82").Value = 2624.5557
$ws.Range("M82").Value = -1168.4
$ws.Range("N82").Value = -3346.5557
$ws.Range("H85").Value = 2233.4285
$ws.Range("I85").Value = 1529.4
$ws.Range("J85").Value = 2624.5557
$ws.Range("K85").Value = 1529.4
$ws.Range("L85").Value = 2624.5557
$ws.Range("M85").Value = -281.4000000000001
$ws.Range("N85").Value = -5120.5557
$ws.Range("H93").Value = 609098.75
$ws.Range("I93").Value = 2221
$ws.Range("K93").Value = 2221
$ws.Range("M93").Value = -973
$ws.Range("H126").Value = 6774.385
$ws.Range("I126").Value = 6507.909
$ws.Range("K126").Value = 19523.727
$ws.Range("M126").Value = -17053.727
$ws.Range("H132").Value = 1681.3214
$ws.Range("J132").Value = 3299
$ws.Range("L132").Value = 9897
$ws.Range("N132").Value = -14957

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 24632.666
$ws.Range("J45").Value = 24632.666
$ws.Range("L45").Value = 24632.666
$ws.Range("N45").Value = -25614.666
$ws.Range("H92").Value = 24500
$ws.Range("J92").Value = 24500
$ws.Range("L92").Value = 24500
$ws.Range("N92").Value = -29492
$ws.Range("H107").Value = 1671.3077
$ws.Range("I107").Value = 1322.7
$ws.Range("J107").Value = 2833.3333
$ws.Range("K107").Value = 3968.1
$ws.Range("L107").Value = 8499.999899999999
$ws.Range("M107").Value = -2048.1
$ws.Range("N107").Value = -12339.9999
$ws.Range("H136").Value = 5105.391
$ws.Range("I136").Value = 3261.9412
$ws.Range("J136").Value = 10328.5
$ws.Range("K136").Value = 9785.8236
$ws.Range("M136").Value = -7235.8236
